$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.378.26"
$ws.Range("E2").Value = "  -2.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.742.52"
$ws.Range("E3").Value = "  -3.29%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.91"
$ws.Range("E5").Value = "  -4.59%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4192"
$ws.Range("E7").Value = "  -8.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3577"
$ws.Range("E8").Value = "  -3.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.45"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07414"
$ws.Range("E10").Value = "  -2.25%  "
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.47"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.078"
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("E15").Value = "  -3.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.738.49"
$ws.Range("E16").Value = "  -3.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001066"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.36"
$ws.Range("E18").Value = "  +6.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06040"
$ws.Range("E19").Value = "  -10.16%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.84"
$ws.Range("E21").Value = "  -4.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.098"
$ws.Range("E22").Value = "  -4.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5225"
$ws.Range("E23").Value = "  -5.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.409.01"
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("E25").Value = "  -3.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.342"
$ws.Range("E26").Value = "  -3.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.42"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.62"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.381"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.934.24"
$ws.Range("E30").Value = "  -3.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.73"
$ws.Range("E31").Value = "  -5.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.178"
$ws.Range("E32").Value = "  -6.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.682"
$ws.Range("E33").Value = "  -2.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09120"
$ws.Range("E34").Value = "  -4.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.623"
$ws.Range("E35").Value = "  -10.10%  "
$ws.Range("E36").Value = "  +4.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02287"
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.065"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06048"
$ws.Range("E40").Value = "  -4.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6375"
$ws.Range("E41").Value = "  -3.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.191"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.429"
$ws.Range("E43").Value = "  -5.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.908"
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.69"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5824"
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.30"
$ws.Range("E49").Value = "  -3.69%  "
$ws.Range("E50").Value = "  -5.16%  "
$ws.Range("E51").Value = "  -4.57%  "
